# Clean-up of input tables:
# Fix the truncated/abbreviated German translations in column C
# of the ID_HeatingSystem lookup table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C2").Value = "Fernwärme"
$ws.Range("C3").Value = "Zentralheizung and Blockheizung"
$ws.Range("C4").Value = "Etagenheizung"
$ws.Range("C5").Value = "Einzel-/Mehrofenheizung"
